$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (people interested) counts in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 134
$wsExhibit.Range("F3").Value = 211
$wsExhibit.Range("F4").Value = 3583
$wsExhibit.Range("F5").Value = 372
$wsExhibit.Range("F6").Value = 21

# Sheet "全部类型" (all types) - same underlying rows, mirror the updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 134
$wsAll.Range("F3").Value = 211
$wsAll.Range("F4").Value = 3583
$wsAll.Range("F5").Value = 372
$wsAll.Range("F8").Value = 21
